# Apply the two logical changes captured by the commit:
#
#  1. Every table that used the deck's local "Table_0" style
#     ({415C1DA5-4DB1-40A1-B04E-BC213222C4B7}) is switched to the
#     built-in "No Style, Table Grid" style
#     ({E5DB8DAB-38C5-4766-B9BB-57588A4FEC63}).
#
#  2. The presentation's theme colour scheme is changed from the
#     "Integral" (Red Violet) palette to the standard "Office Theme"
#     palette.

$p = $ppt.ActivePresentation

$oldTableStyle = "{415C1DA5-4DB1-40A1-B04E-BC213222C4B7}"
$newTableStyle = "{E5DB8DAB-38C5-4766-B9BB-57588A4FEC63}"

for ($idx = 1; $idx -le $p.Slides.Count; $idx++) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldTableStyle) {
                $table.ApplyStyle($newTableStyle)
            }
        }
    }
}

# Re-colour the theme (slide master) to the "Office" palette.
# msoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

function ToRgbValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ToRgbValue($officeThemeColors[$i - 1])
}
